# Fruta / hortaliza, semanal
# Insert 5 new weekly price rows for "Palta" (avocado) - variety "Edranol",
# dated 2021-09-03 (serial 44448), ahead of the existing row 1073 block.
# This pushes the existing rows 1073:1100 down to 1078:1105 and grows the
# sheet's used range to A1:T1105.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 5 blank rows right before the current row 1073.
$ws.Rows("1073:1077").Insert()

# Shared/constant columns for this market-product block.
$mercadoId   = 9
$mercado     = "Vega Central Mapocho de Santiago"
$region      = "Metropolitana"
$codreg      = 13
$tipo        = "Fruta"
$productoId  = 100106
$producto    = "Oleaginosos"
$categoriaId = 100106002
$categoria   = "Palta"
$fecha       = 44448
$unidad      = "$/kilo (en caja de 17 kilos)"
$kgUnidad    = 1

# Row 1073: Edranol - Especial - Provincia de San Felipe de Aconcagua
$r = 1073
$ws.Cells.Item($r, 1).Value  = $mercadoId
$ws.Cells.Item($r, 2).Value  = $mercado
$ws.Cells.Item($r, 3).Value  = $region
$ws.Cells.Item($r, 4).Value  = $fecha
$ws.Cells.Item($r, 5).Value  = $codreg
$ws.Cells.Item($r, 6).Value  = $tipo
$ws.Cells.Item($r, 7).Value  = $productoId
$ws.Cells.Item($r, 8).Value  = $producto
$ws.Cells.Item($r, 9).Value  = $categoriaId
$ws.Cells.Item($r, 10).Value = $categoria
$ws.Cells.Item($r, 11).Value = "Edranol"
$ws.Cells.Item($r, 12).Value = "Especial"
$ws.Cells.Item($r, 13).Value = 115
$ws.Cells.Item($r, 14).Value = 2600
$ws.Cells.Item($r, 15).Value = 2700
$ws.Cells.Item($r, 16).Value = 2648
$ws.Cells.Item($r, 17).Value = $unidad
$ws.Cells.Item($r, 18).Value = "Provincia de San Felipe de Aconcagua"
$ws.Cells.Item($r, 19).Value = 2648
$ws.Cells.Item($r, 20).Value = $kgUnidad

# Row 1074: Edranol - Primera - Provincia de San Felipe de Aconcagua
$r = 1074
$ws.Cells.Item($r, 1).Value  = $mercadoId
$ws.Cells.Item($r, 2).Value  = $mercado
$ws.Cells.Item($r, 3).Value  = $region
$ws.Cells.Item($r, 4).Value  = $fecha
$ws.Cells.Item($r, 5).Value  = $codreg
$ws.Cells.Item($r, 6).Value  = $tipo
$ws.Cells.Item($r, 7).Value  = $productoId
$ws.Cells.Item($r, 8).Value  = $producto
$ws.Cells.Item($r, 9).Value  = $categoriaId
$ws.Cells.Item($r, 10).Value = $categoria
$ws.Cells.Item($r, 11).Value = "Edranol"
$ws.Cells.Item($r, 12).Value = "Primera"
$ws.Cells.Item($r, 13).Value = 95
$ws.Cells.Item($r, 14).Value = 2400
$ws.Cells.Item($r, 15).Value = 2500
$ws.Cells.Item($r, 16).Value = 2447
$ws.Cells.Item($r, 17).Value = $unidad
$ws.Cells.Item($r, 18).Value = "Provincia de San Felipe de Aconcagua"
$ws.Cells.Item($r, 19).Value = 2447
$ws.Cells.Item($r, 20).Value = $kgUnidad

# Row 1075: Edranol - 1a nueva(o) - Provincia de Petorca
$r = 1075
$ws.Cells.Item($r, 1).Value  = $mercadoId
$ws.Cells.Item($r, 2).Value  = $mercado
$ws.Cells.Item($r, 3).Value  = $region
$ws.Cells.Item($r, 4).Value  = $fecha
$ws.Cells.Item($r, 5).Value  = $codreg
$ws.Cells.Item($r, 6).Value  = $tipo
$ws.Cells.Item($r, 7).Value  = $productoId
$ws.Cells.Item($r, 8).Value  = $producto
$ws.Cells.Item($r, 9).Value  = $categoriaId
$ws.Cells.Item($r, 10).Value = $categoria
$ws.Cells.Item($r, 11).Value = "Edranol"
$ws.Cells.Item($r, 12).Value = "1a nueva(o)"
$ws.Cells.Item($r, 13).Value = 145
$ws.Cells.Item($r, 14).Value = 2300
$ws.Cells.Item($r, 15).Value = 2400
$ws.Cells.Item($r, 16).Value = 2352
$ws.Cells.Item($r, 17).Value = $unidad
$ws.Cells.Item($r, 18).Value = "Provincia de Petorca"
$ws.Cells.Item($r, 19).Value = 2352
$ws.Cells.Item($r, 20).Value = $kgUnidad

# Row 1076: Edranol - 2a nueva(o) - Provincia de Petorca
$r = 1076
$ws.Cells.Item($r, 1).Value  = $mercadoId
$ws.Cells.Item($r, 2).Value  = $mercado
$ws.Cells.Item($r, 3).Value  = $region
$ws.Cells.Item($r, 4).Value  = $fecha
$ws.Cells.Item($r, 5).Value  = $codreg
$ws.Cells.Item($r, 6).Value  = $tipo
$ws.Cells.Item($r, 7).Value  = $productoId
$ws.Cells.Item($r, 8).Value  = $producto
$ws.Cells.Item($r, 9).Value  = $categoriaId
$ws.Cells.Item($r, 10).Value = $categoria
$ws.Cells.Item($r, 11).Value = "Edranol"
$ws.Cells.Item($r, 12).Value = "2a nueva(o)"
$ws.Cells.Item($r, 13).Value = 120
$ws.Cells.Item($r, 14).Value = 2100
$ws.Cells.Item($r, 15).Value = 2200
$ws.Cells.Item($r, 16).Value = 2150
$ws.Cells.Item($r, 17).Value = $unidad
$ws.Cells.Item($r, 18).Value = "Provincia de Petorca"
$ws.Cells.Item($r, 19).Value = 2150
$ws.Cells.Item($r, 20).Value = $kgUnidad

# Row 1077: Edranol - Especial nueva (o) - Provincia de Petorca
$r = 1077
$ws.Cells.Item($r, 1).Value  = $mercadoId
$ws.Cells.Item($r, 2).Value  = $mercado
$ws.Cells.Item($r, 3).Value  = $region
$ws.Cells.Item($r, 4).Value  = $fecha
$ws.Cells.Item($r, 5).Value  = $codreg
$ws.Cells.Item($r, 6).Value  = $tipo
$ws.Cells.Item($r, 7).Value  = $productoId
$ws.Cells.Item($r, 8).Value  = $producto
$ws.Cells.Item($r, 9).Value  = $categoriaId
$ws.Cells.Item($r, 10).Value = $categoria
$ws.Cells.Item($r, 11).Value = "Edranol"
$ws.Cells.Item($r, 12).Value = "Especial nueva (o)"
$ws.Cells.Item($r, 13).Value = 130
$ws.Cells.Item($r, 14).Value = 2500
$ws.Cells.Item($r, 15).Value = 2600
$ws.Cells.Item($r, 16).Value = 2546
$ws.Cells.Item($r, 17).Value = $unidad
$ws.Cells.Item($r, 18).Value = "Provincia de Petorca"
$ws.Cells.Item($r, 19).Value = 2546
$ws.Cells.Item($r, 20).Value = $kgUnidad
